$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$ws.Range("B2").Value = 0.0272104810476839
$ws.Range("C2").Value = 0.381361992469342
$ws.Range("D2").Value = 0.3491091507168677
$ws.Range("E2").Value = 0.5908545935480808
$ws.Range("F2").Value = 0.6125082332967703

$ws.Range("B3").Value = 0.1323802532879451
$ws.Range("C3").Value = 0.4216107469797357
$ws.Range("D3").Value = 0.4664420853738647
$ws.Range("E3").Value = 0.6829656546078029
$ws.Range("F3").Value = 0.6973717445328026

$ws.Range("B4").Value = 0.1226621314840486
$ws.Range("C4").Value = 0.511162025793323
$ws.Range("D4").Value = 0.6330782294977604
$ws.Range("E4").Value = 0.7956621327534448
$ws.Range("F4").Value = 0.8211071667389434

$ws.Range("B5").Value = 0.05926794261926783
$ws.Range("C5").Value = 0.5065409700594369
$ws.Range("D5").Value = 0.5385343766439102
$ws.Range("E5").Value = 0.7338490148824282
$ws.Range("F5").Value = 0.7671530853641586

$ws.Range("B6").Value = -0.01880526311880431
$ws.Range("C6").Value = 0.4251950670653798
$ws.Range("D6").Value = 0.3376671520720453
$ws.Range("E6").Value = 0.5810913457211743
$ws.Range("F6").Value = 0.6122032289209177

$ws.Range("B7").Value = -0.04878362568756089
$ws.Range("C7").Value = 0.4489945189615317
$ws.Range("D7").Value = 0.3809843224169152
$ws.Range("E7").Value = 0.6172392748496447
$ws.Range("F7").Value = 0.6526331590693982

$ws.Range("B8").Value = -0.2380091527088908
$ws.Range("C8").Value = 0.3649913975980823
$ws.Range("D8").Value = 0.267128208091335
$ws.Range("E8").Value = 0.5168444718591222
$ws.Range("F8").Value = 0.5025692206868195

$ws.Range("B9").Value = -0.3029994029988738
$ws.Range("C9").Value = 0.4096660696655476
$ws.Range("D9").Value = 0.2533326875858977
$ws.Range("E9").Value = 0.5033216541992781
$ws.Range("F9").Value = 0.4922256332743508

$ws.Range("B10").Value = -0.4690015876295774
$ws.Range("C10").Value = 0.4690015876295774
$ws.Range("D10").Value = 0.2199624891990642
$ws.Range("E10").Value = 0.4690015876295774

